$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.008", "327.26") are not converted to actual numbers,
# matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.620.60"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "1.887.53"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "327.26"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("E7").Value = "  -1.38%  "

$ws.Range("D8").Value = "0.3872"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").Value = "46.83"
$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("D10").Value = "0.07879"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +2.81%  "

$ws.Range("E12").Value = "  -2.69%  "

$ws.Range("D13").Value = "1.894.63"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").Value = "7.089"
$ws.Range("E14").Value = "  +1.13%  "

$ws.Range("D15").Value = "5.721"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").Value = "0.06958"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "87.58"
$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").Value = "0.00001005"
$ws.Range("E19").Value = "  -0.69%  "

$ws.Range("D20").Value = "17.24"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "1.007"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "28.636.52"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").Value = "5.345"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("D25").Value = "2.133.86"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").Value = "2.058"
$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("D27").Value = "154.75"
$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").Value = "5.884"
$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("D30").Value = "1.961"
$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("D31").Value = "118.56"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").Value = "0.09360"
$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("D33").Value = "0.9272"
$ws.Range("E33").Value = "  -1.52%  "

$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("D35").Value = "1.343"
$ws.Range("E35").Value = "  -0.74%  "

$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("D37").Value = "0.05789"
$ws.Range("E37").Value = "  -2.25%  "

$ws.Range("D38").Value = "8.037"
$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("D39").Value = "1.157"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").Value = "0.02075"
$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("D41").Value = "0.5699"

$ws.Range("D42").Value = "0.1800"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").Value = "9.788"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("D46").Value = "0.07154"
$ws.Range("E46").Value = "  -1.56%  "

$ws.Range("D47").Value = "2.174"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("D48").Value = "1.845"
$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("D49").Value = "1.119"
$ws.Range("E49").Value = "  -2.79%  "

$ws.Range("D50").Value = "112.77"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("D51").Value = "2.489"
$ws.Range("E51").Value = "  +4.84%  "

# Row 44/45: Decentraland and EnergySwap swapped positions, with updated values
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "11.81"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5374"
$ws.Range("E45").Value = "  +0.55%  "

# Restore default styling on column D (undo transient NumberFormat change)
$ws.Range("D2:D51").Style = "Normal"
